$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new column after "Total Room" (I) so a fresh "Total Room" column
#    lands at J, and everything from the old J ("Gender") through Q ("Resv Status")
#    shifts right to K..R.
$ws.Range("J1").EntireColumn.Insert()

# 2) The inserted column copied its formatting from its left neighbour (the
#    currency-formatted "Total Rate" column). Re-stamp J2 with the plain
#    General-number look used by the other numeric columns (e.g. E2) before
#    filling in the value, so it doesn't keep the currency format.
$ws.Range("E2").Copy()
$ws.Range("J2").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# 3) Rename / retarget the header text.
$ws.Range("B1").Value = "Nationality"
$ws.Range("H1").Value = "Daily Rate"
$ws.Range("I1").Value = "Total Rate"
$ws.Range("J1").Value = "Total Room"

# 4) Populate the new "Total Room" data cell (same amount as the Total Rate).
$ws.Range("J2").Value = 7500

# 5) The "Nationality" column reads better left aligned than centered.
$ws.Range("B1:B2").HorizontalAlignment = -4131

# 6) Bump the small Arial font used across the booking table from 5pt to 6pt
#    (skip the untouched blank "plain" cells that never used that font).
foreach ($addr in @("A1:R1","A2:A2","C2:J2","L2:L2","N2:O2","Q2:R2")) {
    $ws.Range($addr).Font.Size = 6
}

Write-Host "done"
